# Refresh the "cryptos" price/volume snapshot (GitHub Actions scheduled update).
# Source diff only touches column D (Price) and column E (Volume(1h)) text values
# for data rows 2-51; columns A-C (index / coin name / link) are untouched, as are
# rows where only the Volume(1h) figure moved (Price repeated from the prior run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.464.98"
$ws.Range("E2").Value = "  -5.35%  "
$ws.Range("D3").Value = "1.838.27"
$ws.Range("E3").Value = "  -4.29%  "
$ws.Range("E4").Value = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.46"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -3.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4218"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -7.91%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3638"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -4.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.15"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -3.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07215"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -6.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9023"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -7.69%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.56"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -9.05%  "
$ws.Range("D13").Value = "1.818.20"
$ws.Range("E13").Value = "  -5.41%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.578"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -5.25%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.325"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -6.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06810"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = "  -2.93%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.001"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "77.34"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -8.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000009002"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -5.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9998"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.29"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -8.09%  "
$ws.Range("D22").Value = "27.480.36"
$ws.Range("E22").Value = "  -5.33%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.932"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -7.77%  "
$ws.Range("E24").Value = "  -4.58%  "
$ws.Range("D25").Value = "2.025.55"
$ws.Range("E25").Value = "  -5.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.029"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.03"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.11"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -4.40%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.225"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -6.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "110.61"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -5.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.658"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -9.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08847"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7745"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -9.47%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.499"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -11.40%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.908"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.067"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -13.78%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05346"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -5.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.081"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -5.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01928"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.949"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -4.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.830"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5051"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -7.91%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1629"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -6.93%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.06612"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -4.42%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.203"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -12.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4711"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -8.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.80"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -4.82%  "
$ws.Range("E49").Value = "  -8.39%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.001"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.626"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -7.30%  "
